$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 0, 19.48425592650926)
    3 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 4.327115817150455)
    4 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 8.656069925401464)
    5 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 5.582307763322248)
    6 = @(0.2881169905109251, 0.3048912486333797, 3.223369029078222, 13.86384647080068, 0, 17.68022373902321)
    7 = @(0.1169995834814548, 0.04103571897497393, 3.223369029078222, 0.5333859586016987, 1, 3.914790290136349)
    8 = @(3.272327238179451, 208501.5462402375, 0.7210945179870265, 13.86384647080068, 1, 208519.4035084645)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
